$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "observations" text for the ce_as_002 / ce_as_003 rows first
$ws.Range("C3").Value = "contains null word (length 1)"
$ws.Range("C4").Value = "contains null word (length 0)"

# Append new gr_as_001..gr_as_006 test rows
$ws.Range("A26").Value = "gr_as_001"

$ws.Range("A27").Value = "gr_as_002"
$ws.Range("A28").Value = "gr_as_003"

$ws.Range("B27").Value = "no"
$ws.Range("C27").Value = "double separator"

$ws.Range("B28").Value = "no"
$ws.Range("C28").Value = "no separator"

$ws.Range("A29").Value = "gr_as_004"
$ws.Range("A30").Value = "gr_as_005"
$ws.Range("A31").Value = "gr_as_006"

$ws.Range("B29").Value = "no"
$ws.Range("C29").Value = "separator different than ."

# Last: ce_as_004's observation changes to the new "row with abnormal length" text
$ws.Range("C5").Value = "row with abnormal length"

# Scroll the sheet view so row 4 is at the top, matching the author's saved view
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
